# "progress on section 6" - Ultimate Excel Programmer (VBA) video listing tracker
#
# 1) Rename the active sheet ("Unity Course Video Listing" ->
#    "Excel VBA Course Video Listing") and repoint the Print_Titles
#    defined name at the new sheet name.
# 2) Record completion dates (columns G/H) for a run of Section 6 videos
#    that were watched; everything else (F1/H1 totals, the per-day I/J
#    roll-ups, row spans, sheet dimension, shared-string table growth) is
#    a formula / packaging side effect that recalculates automatically.
# 3) Add a note in K53 explaining why the section is paused.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Rename sheet + fix up the Print_Titles defined name -----------------
$oldName = $ws.Name
$newName = "Excel VBA Course Video Listing"
$ws.Name = $newName

foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Titles*") {
        $n.RefersTo = "='" + $newName + "'!`$3:`$3"
    }
}

# --- 2) Update watched/completed dates in columns G (watched) and H (completed) ---
# Rows 55-60: watched+completed on 42871 (2017-05-15)
for ($r = 55; $r -le 60; $r++) {
    $ws.Cells.Item($r, 7).Value = 42871   # G
    $ws.Cells.Item($r, 8).Value = 42871   # H
}

# Rows 61-68: watched+completed on 42896 (2017-06-09)
for ($r = 61; $r -le 68; $r++) {
    $ws.Cells.Item($r, 7).Value = 42896   # G
    $ws.Cells.Item($r, 8).Value = 42896   # H
}

# Rows 69-76: completed on 42896 (G left as-is)
for ($r = 69; $r -le 76; $r++) {
    $ws.Cells.Item($r, 8).Value = 42896   # H
}

# Rows 78-101: completed on 42897 (2017-06-10)
for ($r = 78; $r -le 101; $r++) {
    $ws.Cells.Item($r, 8).Value = 42897   # H
}

# Rows 103-116: completed on 42897 (2017-06-10)
for ($r = 103; $r -le 116; $r++) {
    $ws.Cells.Item($r, 8).Value = 42897   # H
}

# Rows 118-134: completed on 42898 (2017-06-11)
for ($r = 118; $r -le 134; $r++) {
    $ws.Cells.Item($r, 8).Value = 42898   # H
}

# Rows 136-173: completed on 42898 (2017-06-11)
for ($r = 136; $r -le 173; $r++) {
    $ws.Cells.Item($r, 8).Value = 42898   # H
}

# Rows 175-191: completed on 42899 (2017-06-12)
for ($r = 175; $r -le 191; $r++) {
    $ws.Cells.Item($r, 8).Value = 42899   # H
}

# Rows 193-195: completed on 42899 (2017-06-12)
for ($r = 193; $r -le 195; $r++) {
    $ws.Cells.Item($r, 8).Value = 42899   # H
}

# --- 3) Note why the section is paused ---------------------------------------
$ws.Range("K53").Value = "on hold until after Unity exam"

# --- Leave the selection on G69, matching the author's spot in the sheet ----
$ws.Range("G69").Select()
